$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "Gen" -> "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# Update column A (row 2..14): generation counts -> fraction-of-FES values
$newA = @(0, 0.001, 0.01, 0.1, 0.2, 0.3, 0.4, 0.5, 0.6, 0.7, 0.8, 0.9, 1)
for ($i = 0; $i -lt $newA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newA[$i]
}

# New recomputed "Mean" values (replacing the old Run 50 + Mean columns)
$newMean = @(13.92583074, 12.6104404, 9.97376794, 6.23819155, 5.1548993, 4.51799287, 4.00897376, 3.67560886, 3.30055844, 3.07075649, 2.80947207, 2.57491379, 2.37759439)
for ($i = 0; $i -lt $newMean.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 52).Value = $newMean[$i]
}

# Rename the last-column header from "Run 50" to "Mean"
$ws.Cells.Item(1, 52).Value = "Mean"

# Delete the now-obsolete last column (BA, "Mean" / old Run50 data column shifted)
$ws.Columns.Item(53).Delete()
